$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E7").Value = 15.176
$ws.Range("A8").Value = -22.34230000000001
$ws.Range("A10").Value = -21.87989999999998
$ws.Range("A12").Value = -21.4889
$ws.Range("E14").Value = 16.7854
$ws.Range("E15").Value = 15.9731
$ws.Range("A18").Value = -22.07610000000001
$ws.Range("E18").Value = 17.68230000000002
$ws.Range("E20").Value = 15.894
$ws.Range("A25").Value = -21.56209999999998
$ws.Range("E29").Value = 17.09580000000001
$ws.Range("E30").Value = 15.6206
$ws.Range("E31").Value = 16.04550000000001
$ws.Range("E35").Value = 15.767
$ws.Range("A37").Value = -19.57039999999999
$ws.Range("E40").Value = 17.1032
$ws.Range("E44").Value = 16.50939999999999
$ws.Range("E50").Value = 16.5152
$ws.Range("E54").Value = 16.6226
$ws.Range("A55").Value = -22.31070000000001
$ws.Range("A68").Value = -21.4838
$ws.Range("E68").Value = 17.11260000000001
$ws.Range("E76").Value = 16.16379999999998
$ws.Range("A77").Value = -20.23729999999999
$ws.Range("A78").Value = -19.74439999999998
$ws.Range("A79").Value = -20.29189999999999
$ws.Range("A80").Value = -19.3995
$ws.Range("A81").Value = -21.7756
$ws.Range("A82").Value = -22.1
$ws.Range("A84").Value = -22.0652
$ws.Range("E87").Value = 16.13839999999999
$ws.Range("E88").Value = 16.3741
$ws.Range("E92").Value = 18.50910000000002
$ws.Range("E96").Value = 16.1549
$ws.Range("E98").Value = 15.5709
$ws.Range("A101").Value = -21.32629999999999
$ws.Range("E101").Value = 16.8127
$ws.Range("A102").Value = -19.55469999999999
$ws.Range("E102").Value = 16.5487
